# "Added first version of clicking based on excel data"
# - rename the sheet from the Finnish default "Taul1" to "Sheet1"
# - replace the two header/label cells with their new English text
# - move the active selection and bump the zoom level, matching the
#   view state captured when the sheet was last saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "                          Ostot                      "

$excel.ActiveWindow.Zoom = 235
$ws.Range("H6").Select()
